# Apply targeted value updates to the "Jay Test 2 - LSMeans" worksheet
# per commit: "Working on functions to do the combined analysis with the 2020/2021 data"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Jay Test 2 - LSMeans")

$ws.Range("Y5").Value = 63
$ws.Range("AA5").Value = 54.9

$ws.Range("T9").Value = 33.8
$ws.Range("AB9").Value = 2266

$ws.Range("T10").Value = 31.8

$ws.Range("T11").Value = 32.8

$ws.Range("T12").Value = 33.3

$ws.Range("T14").Value = 37.8

$ws.Range("T16").Value = 34.8

$ws.Range("T18").Value = 31.8

$ws.Range("T19").Value = 31.3

$ws.Range("P20").Value = 67.1

$ws.Range("T21").Value = 36.8

$ws.Range("R22").Value = 57.5
$ws.Range("T22").Value = 35.8
$ws.Range("AA22").Value = 57.1
